# Apply updated crypto price/volume figures to the sheet.
# Values are entered with a leading apostrophe so Excel stores them
# as literal text (matching the original inlineStr cell type) instead
# of auto-converting number-like strings (e.g. "206.86") into numerics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$textPrefix = "'"

$c = $ws.Range('D2')
$c.Value = $textPrefix + '27.685.07'
$c.ClearFormats()

$c = $ws.Range('E2')
$c.Value = $textPrefix + '  -0.45%  '
$c.ClearFormats()

$c = $ws.Range('D3')
$c.Value = $textPrefix + '1.583.97'
$c.ClearFormats()

$c = $ws.Range('E3')
$c.Value = $textPrefix + '  -2.56%  '
$c.ClearFormats()

$c = $ws.Range('E4')
$c.Value = $textPrefix + '  +0.76%  '
$c.ClearFormats()

$c = $ws.Range('D5')
$c.Value = $textPrefix + '206.86'
$c.ClearFormats()

$c = $ws.Range('E5')
$c.Value = $textPrefix + '  -1.92%  '
$c.ClearFormats()

$c = $ws.Range('D6')
$c.Value = $textPrefix + '0.501'
$c.ClearFormats()

$c = $ws.Range('E6')
$c.Value = $textPrefix + '  -3.36%  '
$c.ClearFormats()

$c = $ws.Range('E7')
$c.Value = $textPrefix + '  +0.74%  '
$c.ClearFormats()

$c = $ws.Range('D8')
$c.Value = $textPrefix + '22.12'
$c.ClearFormats()

$c = $ws.Range('E8')
$c.Value = $textPrefix + '  -4.68%  '
$c.ClearFormats()

$c = $ws.Range('E9')
$c.Value = $textPrefix + '  -1.59%  '
$c.ClearFormats()

$c = $ws.Range('E10')
$c.Value = $textPrefix + '  -2.82%  '
$c.ClearFormats()

$c = $ws.Range('D11')
$c.Value = $textPrefix + '0.0866'
$c.ClearFormats()

$c = $ws.Range('E11')
$c.Value = $textPrefix + '  -1.40%  '
$c.ClearFormats()

$c = $ws.Range('E12')
$c.Value = $textPrefix + '  -2.55%  '
$c.ClearFormats()

$c = $ws.Range('D13')
$c.Value = $textPrefix + '1.605.69'
$c.ClearFormats()

$c = $ws.Range('E13')
$c.Value = $textPrefix + '  -1.14%  '
$c.ClearFormats()

$c = $ws.Range('E14')
$c.Value = $textPrefix + '  -4.20%  '
$c.ClearFormats()

$c = $ws.Range('D15')
$c.Value = $textPrefix + '0.527'
$c.ClearFormats()

$c = $ws.Range('E15')
$c.Value = $textPrefix + '  -5.06%  '
$c.ClearFormats()

$c = $ws.Range('D17')
$c.Value = $textPrefix + '27.630.54'
$c.ClearFormats()

$c = $ws.Range('E17')
$c.Value = $textPrefix + '  -0.71%  '
$c.ClearFormats()

$c = $ws.Range('D18')
$c.Value = $textPrefix + '219.38'
$c.ClearFormats()

$c = $ws.Range('E18')
$c.Value = $textPrefix + '  -3.94%  '
$c.ClearFormats()

$c = $ws.Range('D19')
$c.Value = $textPrefix + '0.0₃0694'
$c.ClearFormats()

$c = $ws.Range('E19')
$c.Value = $textPrefix + '  -3.21%  '
$c.ClearFormats()

$c = $ws.Range('E20')
$c.Value = $textPrefix + '  -3.59%  '
$c.ClearFormats()

$c = $ws.Range('E21')
$c.Value = $textPrefix + '  +0.85%  '
$c.ClearFormats()

$c = $ws.Range('D22')
$c.Value = $textPrefix + '4.12'
$c.ClearFormats()

$c = $ws.Range('E22')
$c.Value = $textPrefix + '  -4.54%  '
$c.ClearFormats()

$c = $ws.Range('D23')
$c.Value = $textPrefix + '9.53'
$c.ClearFormats()

$c = $ws.Range('E23')
$c.Value = $textPrefix + '  -3.77%  '
$c.ClearFormats()

$c = $ws.Range('E24')
$c.Value = $textPrefix + '  -3.63%  '
$c.ClearFormats()

$c = $ws.Range('D25')
$c.Value = $textPrefix + '153.62'
$c.ClearFormats()

$c = $ws.Range('E25')
$c.Value = $textPrefix + '  -0.88%  '
$c.ClearFormats()

$c = $ws.Range('E26')
$c.Value = $textPrefix + '  -0.93%  '
$c.ClearFormats()

$c = $ws.Range('E27')
$c.Value = $textPrefix + '  +0.83%  '
$c.ClearFormats()

$c = $ws.Range('D28')
$c.Value = $textPrefix + '15.08'
$c.ClearFormats()

$c = $ws.Range('E28')
$c.Value = $textPrefix + '  -2.16%  '
$c.ClearFormats()

$c = $ws.Range('E29')
$c.Value = $textPrefix + '  -4.53%  '
$c.ClearFormats()

$c = $ws.Range('E30')
$c.Value = $textPrefix + '  -1.90%  '
$c.ClearFormats()

$c = $ws.Range('D31')
$c.Value = $textPrefix + '0.0465'
$c.ClearFormats()

$c = $ws.Range('E31')
$c.Value = $textPrefix + '  -2.97%  '
$c.ClearFormats()

$c = $ws.Range('E32')
$c.Value = $textPrefix + '  -5.65%  '
$c.ClearFormats()

$c = $ws.Range('D33')
$c.Value = $textPrefix + '1.359.78'
$c.ClearFormats()

$c = $ws.Range('E33')
$c.Value = $textPrefix + '  -2.93%  '
$c.ClearFormats()

$c = $ws.Range('E34')
$c.Value = $textPrefix + '  -4.96%  '
$c.ClearFormats()

$c = $ws.Range('D35')
$c.Value = $textPrefix + '1.52'
$c.ClearFormats()

$c = $ws.Range('E35')
$c.Value = $textPrefix + '  -4.49%  '
$c.ClearFormats()

$c = $ws.Range('D36')
$c.Value = $textPrefix + '0.971'
$c.ClearFormats()

$c = $ws.Range('E36')
$c.Value = $textPrefix + '  -3.03%  '
$c.ClearFormats()

$c = $ws.Range('D37')
$c.Value = $textPrefix + '2.31'
$c.ClearFormats()

$c = $ws.Range('E37')
$c.Value = $textPrefix + '  -0.61%  '
$c.ClearFormats()

$c = $ws.Range('E38')
$c.Value = $textPrefix + '  -1.39%  '
$c.ClearFormats()

$c = $ws.Range('D39')
$c.Value = $textPrefix + '0.535'
$c.ClearFormats()

$c = $ws.Range('E39')
$c.Value = $textPrefix + '  -3.26%  '
$c.ClearFormats()

$c = $ws.Range('D40')
$c.Value = $textPrefix + '0.819'
$c.ClearFormats()

$c = $ws.Range('E40')
$c.Value = $textPrefix + '  -3.10%  '
$c.ClearFormats()

$c = $ws.Range('E41')
$c.Value = $textPrefix + '  +0.72%  '
$c.ClearFormats()

$c = $ws.Range('D42')
$c.Value = $textPrefix + '0.967'
$c.ClearFormats()

$c = $ws.Range('E42')
$c.Value = $textPrefix + '  -3.03%  '
$c.ClearFormats()

$c = $ws.Range('D43')
$c.Value = $textPrefix + '63.59'
$c.ClearFormats()

$c = $ws.Range('E43')
$c.Value = $textPrefix + '  -2.98%  '
$c.ClearFormats()

$c = $ws.Range('D44')
$c.Value = $textPrefix + '2.16'
$c.ClearFormats()

$c = $ws.Range('E44')
$c.Value = $textPrefix + '  +2.36%  '
$c.ClearFormats()

$c = $ws.Range('D45')
$c.Value = $textPrefix + '5.18'
$c.ClearFormats()

$c = $ws.Range('E45')
$c.Value = $textPrefix + '  -3.96%  '
$c.ClearFormats()

$c = $ws.Range('E46')
$c.Value = $textPrefix + '  -4.50%  '
$c.ClearFormats()

$c = $ws.Range('D47')
$c.Value = $textPrefix + '1.719.85'
$c.ClearFormats()

$c = $ws.Range('E47')
$c.Value = $textPrefix + '  -2.57%  '
$c.ClearFormats()

$c = $ws.Range('D48')
$c.Value = $textPrefix + '87.98'
$c.ClearFormats()

$c = $ws.Range('E48')
$c.Value = $textPrefix + '  -0.16%  '
$c.ClearFormats()

$c = $ws.Range('E49')
$c.Value = $textPrefix + '  +11.57%  '
$c.ClearFormats()

$c = $ws.Range('D50')
$c.Value = $textPrefix + '0.0967'
$c.ClearFormats()

$c = $ws.Range('E50')
$c.Value = $textPrefix + '  -4.11%  '
$c.ClearFormats()

$c = $ws.Range('E51')
$c.Value = $textPrefix + '  -1.10%  '
$c.ClearFormats()

